$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("N7").Value = 10
$ws.Range("AC7").Value = 10
$ws.Range("AD7").Value = 8.5
$ws.Range("AF7").Value = 81
$ws.Range("AG7").Value = 15
$ws.Range("AH7").Value = 34
